$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(82)
$r = $p.Range
Write-Output ("WordOpenXML: {0}" -f $r.WordOpenXML)
